$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 17): Number of employees / Assets / Turnover ---
$ws.Range("B17").Value = "Number of employees"
$ws.Range("B17").Style = "title"
$ws.Range("C17").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("C17").Style = "title"
$ws.Range("D17").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("D17").Style = "title"

# --- Row 18: Micro ---
$ws.Range("A18").Value = "Micro"
$ws.Range("A18").Style = "Normal"
$ws.Range("B18").Value = "'"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'"
$ws.Range("D18").Style = "Normal"

# --- Row 19: Small ---
$ws.Range("A19").Value = "Small"
$ws.Range("A19").Style = "Normal"
$ws.Range("B19").Value = "'"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'"
$ws.Range("D19").Style = "Normal"

# --- Row 20: Medium (overwrites old "MFA" cell) ---
$ws.Range("A20").Value = "Medium"
$ws.Range("A20").Style = "Normal"
$ws.Range("B20").Value = "<300"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'"
$ws.Range("D20").Style = "Normal"

# --- Row 21: Large (overwrites old source-citation cell) ---
$ws.Range("A21").Value = "Large"
$ws.Range("A21").Style = "Normal"
$ws.Range("B21").Value = ">=300"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'"
$ws.Range("D21").Style = "Normal"

# --- Rows 26-27: relocate the original "MFA" source citation ---
$ws.Range("A26").Value = "MFA"
$ws.Range("A26").Style = "title"
$ws.Range("A27").Value = "Ministry of Foreign Affaris (MFA), ""SMALL TO MEDIUM ENTERPRISE PAPERS"", N/S, p. 3. Available at http://mfa.gov.af/content/files/SME%20PAPER.pdf"
$ws.Range("A27").Style = "source"
